$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "2nd PC" result table (columns K:M, rows 6-11)
$ws.Range("K6").Value = 32
$ws.Range("L6").Value = 30
$ws.Range("M6").Value = 30

$ws.Range("K7").Value = 27
$ws.Range("L7").Value = 28
$ws.Range("M7").Value = 27

$ws.Range("K8").Value = 25
$ws.Range("L8").Value = 25
$ws.Range("M8").Value = 25

$ws.Range("K9").Value = 23
$ws.Range("L9").Value = 22
$ws.Range("M9").Value = 22

$ws.Range("K10").Value = 21
$ws.Range("L10").Value = 22
$ws.Range("M10").Value = 21

$ws.Range("K11").Value = 22
$ws.Range("L11").Value = 22
$ws.Range("M11").Value = 22

# Clear the obsolete "On second with 2 threads / 16 sec / With 4: 12 sec" note
# (A11:B11 keep their style but lose their value; C11:F11 are cleared outright)
$ws.Range("A11:B11").ClearContents()
$ws.Range("C11:F11").ClearContents()

# Clear the obsolete "On third almost the same as on second" note (row 12)
$ws.Range("A12:E12").ClearContents()

# Move the selection back to A11 and scroll the view to the top
$ws.Range("A1").Select()
$ws.Range("A11").Select()
